# Add "2022-Q4" quarterly data to the BAS-巴斯夫欧洲公司 workbook.
#
# Net effect (per the source diff):
#   1. A brand-new worksheet named "2022-Q4" is inserted right after the
#      "总计" (summary) sheet and right before the existing "2022-Q3" sheet.
#      It carries one fund-holding data row (same shape/style as the other
#      quarterly sheets).
#   2. The "总计" summary sheet gets a new data row inserted right below its
#      header (pushing the existing quarter rows down by one), recording the
#      new quarter's holding count / market value; the leading numeric index
#      column is renumbered sequentially.
#   3. Every other quarterly sheet (2022-Q3 ... 2020-Q4) keeps its own name
#      and its own data untouched - they simply shift one tab position to
#      the right to make room for the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q3" sheet (so the new sheet inherits the
#    exact same layout/formatting) and place the copy right after "总计".
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $summarySheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Overwrite the fund figures with the new quarter's numbers. Column D-G hold
# text-formatted numbers in this workbook (leading apostrophe keeps them as
# text instead of auto-converting to a number), column H is a plain number.
$q4Sheet.Range("D2").Value = "'4.76"
$q4Sheet.Range("E2").Value = "'92.90"
$q4Sheet.Range("F2").Value = "'3.31"
$q4Sheet.Range("G2").Value = "'0.1576"
$q4Sheet.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: shift the 8 existing data rows (rows 2-9) down
#    into rows 3-10 (carrying their formatting along), then write the new
#    "2022-Q4" row into row 2 and renumber the index column (A).
# ---------------------------------------------------------------------
$summarySheet.Range("A2:D9").Copy($summarySheet.Range("A3:D10"))

$summarySheet.Range("B2").Value = "2022-Q4"
$summarySheet.Range("C2").Value = 1
$summarySheet.Range("D2").Value = 0.16

for ($row = 2; $row -le 10; $row++) {
    $summarySheet.Cells.Item($row, 1).Value = $row - 2
}

# ---------------------------------------------------------------------
# 3. Restore the originally active tab ("2020-Q4", the last sheet) - adding
#    the new sheet shifts which tab is marked selected, so re-activate it.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
